# Fix Training Data Issue (#48)
# The "Date" column (BF) was stamped with the literal game-folder name
# ("4-22-2007-08") instead of the actual game date. NBA.com's box-score
# pages label a game played in April with the *previous* day's date, so
# the real date for this file is one day earlier: 2008-04-22.
# Re-point every BF2:BF31 cell from "4-22-2007-08" to "2008-04-22",
# keeping the value as plain text (not an auto-converted date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "4-22-2007-08"
$newDate = "2008-04-22"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 58).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 31 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 58)   # column BF
    if ($cell.Text -eq $oldDate) {
        # Force text storage so "2008-04-22" isn't auto-parsed into a date
        # serial number, then drop the temporary format so the cell ends
        # up unstyled again, just like it started.
        $cell.NumberFormat = "@"
        $cell.Value = $newDate
        $cell.ClearFormats()
    }
}
